$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "04:00 PM"
$ws.Range("I2").Value = "04:00 PM"
$ws.Range("G3").Value = "04:00 PM"
$ws.Range("I3").Value = "04:00 PM"

$ws.Range("I3").Select()
